$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.878.21'

$ws.Range("D3").Value = '1.887.39'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7486'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '

$ws.Range("E7").Value = '  +0.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3117'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.40'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07114'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08533'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.97%  '

$ws.Range("D13").Value = '1.892.40'
$ws.Range("E13").Value = '  -3.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.359'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.136'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.98%  '

$ws.Range("D17").Value = '29.959.58'
$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("E18").Value = '  -2.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007792'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.76%  '

$ws.Range("D21").Value = '2.155.13'
$ws.Range("E21").Value = '  +3.22%  '

$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.996'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1597'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.355'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.026'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.518'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.96%  '

$ws.Range("E31").Value = '  -0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.468'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.092'
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05383'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.235'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7423'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.003'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.713'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.57%  '

$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.771'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.13%  '

$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("D42").Value = '1.103.25'
$ws.Range("E42").Value = '  -3.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.069'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8568'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.002'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.47'
$ws.Range("D47").Style = "Normal"

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.655'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.863'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.058'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.93%  '

$ws.Range("D51").Value = '2.045.28'
$ws.Range("E51").Value = '  +3.10%  '

